$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.305.06"
$ws.Range("E2").Value = "  -2.69%  "
$ws.Range("D3").Value = "'3.306.15"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'559.09"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "'141.94"
$ws.Range("E6").Value = "  -4.17%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'3.306.62"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").Value = "'0.473"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "'0.118"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "'3.878.27"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "'26.80"
$ws.Range("E15").Value = "  -5.60%  "
$ws.Range("D16").Value = "'3.305.24"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "'60.333.20"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("D20").Value = "'14.44"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "'8.64"
$ws.Range("E21").Value = "  -3.56%  "
$ws.Range("D22").Value = "'374.64"
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("D23").Value = "'74.21"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  -4.35%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'3.444.43"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("E27").Value = "  -7.49%  "
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "'7.25"
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "'7.63"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("D34").Value = "'22.56"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "'1.27"
$ws.Range("E35").Value = "  -5.12%  "
$ws.Range("D36").Value = "'5.13"
$ws.Range("E36").Value = "  -6.10%  "
$ws.Range("D37").Value = "'166.18"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("E38").Value = "  -5.07%  "
$ws.Range("D39").Value = "'6.72"
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.99"
$ws.Range("E40").Value = "  -10.64%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "'3.339.50"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "'0.0736"
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("D43").Value = "'42.00"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").Value = "'2.374.81"
$ws.Range("E48").Value = "  -6.55%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'6.49"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "'21.28"
$ws.Range("E51").Value = "  -6.41%  "
